$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Instructions")

# Sheet is protected; unprotect before making structural changes
$ws.Unprotect()

# Update version number
$ws.Range("A2").Value = "Version 1.2.3"

# Update text of row 5 (now a different instruction)
$ws.Range("A5").Value = "Please use consecutive rows (no blank rows)."

# Insert a new row at row 6 (pushes existing row 6 and below down by one)
$ws.Rows.Item(6).Insert()

# Set the text for the newly inserted row 6
$ws.Range("A6").Value = "Do not edit the header row of the 'Antibodies' sheet."

# Re-apply sheet protection matching original settings
$ws.Protect()
